$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 484
$ws.Range("I38").Value = 61.909092
$ws.Range("J38").Value = 1644.75
$ws.Range("K38").Value = 185.727276
$ws.Range("L38").Value = 4934.25
$ws.Range("M38").Value = 186.272724
$ws.Range("N38").Value = -5678.25

$ws.Range("H39").Value = 300.23077
$ws.Range("J39").Value = 721.2
$ws.Range("L39").Value = 2163.6
$ws.Range("N39").Value = -2755.6

$ws.Range("H58").Value = 316.16666
$ws.Range("J58").Value = 500
$ws.Range("L58").Value = 1500
$ws.Range("N58").Value = -1800

$ws.Range("H76").Value = 4071.5
$ws.Range("I76").Value = 4071.5
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 4071.5
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -3756.5

$ws.Range("H79").Value = 4071.5
$ws.Range("I79").Value = 4071.5
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 4071.5
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -2979.5

$ws.Range("H135").Value = 1100.32
$ws.Range("I135").Value = 680
$ws.Range("K135").Value = 6120
$ws.Range("M135").Value = -3585

$ws.Range("H137").Value = 1865.2609
$ws.Range("I137").Value = 1666.9445
$ws.Range("J137").Value = 2579.2
$ws.Range("K137").Value = 5000.833500000001
$ws.Range("L137").Value = 7737.599999999999
$ws.Range("M137").Value = -2450.833500000001
$ws.Range("N137").Value = -12837.6

$ws.Range("H138").Value = 3214.1973
$ws.Range("I138").Value = 1541.8077
$ws.Range("J138").Value = 4083.84
$ws.Range("K138").Value = 4625.4231
$ws.Range("L138").Value = 12251.52
$ws.Range("M138").Value = 514.5769
$ws.Range("N138").Value = -22531.52

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 140000
$ws.Range("J23").Value = 140000
$ws.Range("L23").Value = 140000
$ws.Range("N23").Value = -140518

$ws.Range("H61").Value = 5866
$ws.Range("I61").Value = 5700.875
$ws.Range("J61").Value = 6130.2
$ws.Range("K61").Value = 5700.875
$ws.Range("L61").Value = 6130.2
$ws.Range("M61").Value = -5488.875
$ws.Range("N61").Value = -6554.2

$ws.Range("H74").Value = 13762.866
$ws.Range("I74").Value = 1874.8
$ws.Range("J74").Value = 37539
$ws.Range("K74").Value = 1874.8
$ws.Range("L74").Value = 37539
$ws.Range("M74").Value = -1000.8
$ws.Range("N74").Value = -39287

$ws.Range("H77").Value = 13762.866
$ws.Range("I77").Value = 1874.8
$ws.Range("J77").Value = 37539
$ws.Range("K77").Value = 9374
$ws.Range("L77").Value = 187695
$ws.Range("M77").Value = -5006
$ws.Range("N77").Value = -196431

$ws.Range("H97").Value = 1288.6875
$ws.Range("I97").Value = 1386.3572
$ws.Range("J97").Value = 605
$ws.Range("K97").Value = 1386.3572
$ws.Range("L97").Value = 605
$ws.Range("M97").Value = -890.3571999999999
$ws.Range("N97").Value = -1597

$ws.Range("H132").Value = 5547.909
$ws.Range("I132").Value = 4273.75
$ws.Range("J132").Value = 6276
$ws.Range("K132").Value = 12821.25
$ws.Range("L132").Value = 18828
$ws.Range("M132").Value = -10291.25
$ws.Range("N132").Value = -23888

$ws.Range("H136").Value = 5866
$ws.Range("I136").Value = 5700.875
$ws.Range("J136").Value = 6130.2
$ws.Range("K136").Value = 17102.625
$ws.Range("L136").Value = 18390.6
$ws.Range("M136").Value = -14552.625
$ws.Range("N136").Value = -23490.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 849.5357
$ws.Range("I107").Value = 679.8095
$ws.Range("J107").Value = 1358.7142
$ws.Range("K107").Value = 679.8095
$ws.Range("L107").Value = 1358.7142
$ws.Range("M107").Value = 1240.1905
$ws.Range("N107").Value = -5198.7142

$ws.Range("H134").Value = 2409.7856
$ws.Range("I134").Value = 2424.963
$ws.Range("K134").Value = 7274.889000000001
$ws.Range("M134").Value = -4739.889000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1262.4445
$ws.Range("I16").Value = 974.75
$ws.Range("J16").Value = 1492.6
$ws.Range("K16").Value = 974.75
$ws.Range("L16").Value = 1492.6
$ws.Range("M16").Value = -687.75
$ws.Range("N16").Value = -2066.6

$ws.Range("H22").Value = 774.75
$ws.Range("I22").Value = 550
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 550
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = -200
$ws.Range("N22").Value = -1699.5

$ws.Range("H31").Value = 47200.22
$ws.Range("I31").Value = 51326.8
$ws.Range("J31").Value = 19689.666
$ws.Range("K31").Value = 51326.8
$ws.Range("L31").Value = 19689.666
$ws.Range("M31").Value = -51031.8
$ws.Range("N31").Value = -20279.666

$ws.Range("H34").Value = 47200.22
$ws.Range("I34").Value = 51326.8
$ws.Range("J34").Value = 19689.666
$ws.Range("K34").Value = 51326.8
$ws.Range("L34").Value = 19689.666
$ws.Range("M34").Value = -51124.8
$ws.Range("N34").Value = -20093.666

$ws.Range("H107").Value = 789.64703
$ws.Range("I107").Value = 468.1
$ws.Range("J107").Value = 1249
$ws.Range("K107").Value = 468.1
$ws.Range("L107").Value = 1249
$ws.Range("M107").Value = 1451.9
$ws.Range("N107").Value = -5089

$ws.Range("H113").Value = 1262.4445
$ws.Range("I113").Value = 974.75
$ws.Range("J113").Value = 1492.6
$ws.Range("K113").Value = 974.75
$ws.Range("L113").Value = 1492.6
$ws.Range("M113").Value = 1195.25
$ws.Range("N113").Value = -5832.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 110.4
$ws.Range("I15").Value = 17.333334
$ws.Range("K15").Value = 52.000002
$ws.Range("M15").Value = 87.99999800000001

$ws.Range("H17").Value = 31
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 13007.143
$ws.Range("I57").Value = 8027.5
$ws.Range("K57").Value = 8027.5
$ws.Range("M57").Value = -7207.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 17627.75
$ws.Range("I17").Value = 23500.334
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = 23500.334
$ws.Range("L17").Value = 10
$ws.Range("M17").Value = -23330.334
$ws.Range("N17").Value = -350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11043.429
$ws.Range("I81").Value = 34037
$ws.Range("K81").Value = 68074
$ws.Range("M81").Value = -67013

$ws.Range("H84").Value = 11043.429
$ws.Range("I84").Value = 34037
$ws.Range("K84").Value = 340370
$ws.Range("M84").Value = -335066

$ws.Range("H132").Value = 2721.9207
$ws.Range("I132").Value = 2480.224
$ws.Range("K132").Value = 7274.889000000001
$ws.Range("M132").Value = -4910.672
